$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new note in L4 (write first so it claims the earlier shared-string slot)
$ws.Range("L4").Value = "seems fine"

# Update the "rejection rate" note from "low" to "high"
$ws.Range("E6").Value = "rejection rate for f__Lachnospiraceae is high"

# Update the window selection / scroll position
[void]$ws.Range("E6").Select()
